$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32 ("primarykey"): B32 was blank, now set to "id"
$ws.Range("B32").Value = "id"

# Row 33 ("samplelimit"): B33 changes from 5 to 8 and becomes left-aligned
$ws.Range("B33").Value = 8
$ws.Range("B33").HorizontalAlignment = -4131  # xlLeft

# Reflect the recorded session's scrolled viewport and final selection
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B34").Select()
